# Updated cryptos list on Sat Mar 30 20:37:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '69.667.97'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.504.97'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "`'" + '605.14'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").Value = "`'" + '194.93'
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("D7").Value = "`'" + '0.626'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = "`'" + '0.201'
$ws.Range("E9").Value = '  -5.69%  '
$ws.Range("D11").Value = "`'" + '53.46'
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("D12").Value = "`'" + '0.0000300'
$ws.Range("E12").Value = '  -2.11%  '
$ws.Range("D13").Value = "`'" + '9.51'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").Value = '4.065.57'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").Value = "`'" + '594.12'
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").Value = "`'" + '19.11'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '69.846.83'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = "`'" + '12.76'
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("E19").Value = '  +2.16%  '
$ws.Range("D20").Value = '3.508.64'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").Value = "`'" + '0.988'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = "`'" + '18.38'
$ws.Range("E22").Value = '  +6.94%  '
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").Value = "`'" + '4.64'
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("D25").Value = "`'" + '101.65'
$ws.Range("E25").Value = '  -3.51%  '
$ws.Range("E26").Value = '  +4.24%  '
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = "`'" + '9.52'
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").Value = "`'" + '33.22'
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("D31").Value = "`'" + '4.29'
$ws.Range("E31").Value = '  +3.29%  '
$ws.Range("D32").Value = "`'" + '12.41'
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("D35").Value = '0.0₃0821'
$ws.Range("E35").Value = '  +6.21%  '
$ws.Range("D36").Value = '3.732.36'
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("D37").Value = "`'" + '3.09'
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("D38").Value = "`'" + '0.999'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").Value = "`'" + '36.38'
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("D42").Value = "`'" + '483.66'
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("D44").Value = "`'" + '0.0453'
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("E46").Value = '  -3.24%  '
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = "`'" + '8.40'
$ws.Range("E49").Value = '  -4.30%  '
$ws.Range("E50").Value = '  +2.55%  '
$ws.Range("E51").Value = '  +10.27%  '
